$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# R30 rule row: set the "Integer min" value (column C, row 10) from 18 to 1.
$ws.Range("C10").Value = 1
